$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) and E (Volume 1h) to Text format so numeric-looking
# strings (e.g. "69.469.78", "11.40", "0.0000175") are preserved exactly as text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Column D (Price) updates
$ws.Range("D2").Value = '69.469.78'
$ws.Range("D3").Value = '2.509.97'
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").Value = '575.49'
$ws.Range("D6").Value = '166.05'
$ws.Range("D8").Value = '0.516'
$ws.Range("D9").Value = '2.507.80'
$ws.Range("D10").Value = '0.159'
$ws.Range("D12").Value = '0.342'
$ws.Range("D13").Value = '4.84'
$ws.Range("D14").Value = '2.967.94'
$ws.Range("D15").Value = '69.395.76'
$ws.Range("D16").Value = '0.0000175'
$ws.Range("D17").Value = '24.72'
$ws.Range("D18").Value = '2.504.03'
$ws.Range("D19").Value = '11.40'
$ws.Range("D20").Value = '7.77'
$ws.Range("D21").Value = '348.31'
$ws.Range("D22").Value = '3.93'
$ws.Range("D23").Value = '1.96'
$ws.Range("D25").Value = '68.34'
$ws.Range("D26").Value = '3.99'
$ws.Range("D27").Value = '8.92'
$ws.Range("D28").Value = '2.637.86'
$ws.Range("D29").Value = '0.997'
$ws.Range("D30").Value = '0.0₃0898'
$ws.Range("D31").Value = '7.87'
$ws.Range("D32").Value = '472.64'
$ws.Range("D33").Value = '1.27'
$ws.Range("D36").Value = '0.116'
$ws.Range("D37").Value = '153.53'
$ws.Range("D38").Value = '18.94'
$ws.Range("D39").Value = '18.45'
$ws.Range("D41").Value = '4.73'
$ws.Range("D42").Value = '0.316'
$ws.Range("D46").Value = '38.08'
$ws.Range("D47").Value = '143.03'
$ws.Range("D48").Value = '0.527'
$ws.Range("D49").Value = '3.52'

# Column E (Volume 1h) updates
$ws.Range("E2").Value = '  -3.87%  '
$ws.Range("E3").Value = '  -5.03%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E5").Value = '  -2.29%  '
$ws.Range("E6").Value = '  -4.73%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.56%  '
$ws.Range("E9").Value = '  -5.08%  '
$ws.Range("E10").Value = '  -7.24%  '
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("E12").Value = '  -3.81%  '
$ws.Range("E13").Value = '  -2.07%  '
$ws.Range("E14").Value = '  -5.24%  '
$ws.Range("E15").Value = '  -3.84%  '
$ws.Range("E16").Value = '  -5.97%  '
$ws.Range("E17").Value = '  -4.49%  '
$ws.Range("E18").Value = '  -4.21%  '
$ws.Range("E19").Value = '  -6.78%  '
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("E21").Value = '  -6.74%  '
$ws.Range("E22").Value = '  -5.30%  '
$ws.Range("E23").Value = '  -4.99%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").Value = '  -3.85%  '
$ws.Range("E26").Value = '  -6.58%  '
$ws.Range("E27").Value = '  -7.43%  '
$ws.Range("E28").Value = '  -5.06%  '
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("E30").Value = '  -5.88%  '
$ws.Range("E31").Value = '  -1.81%  '
$ws.Range("E32").Value = '  -4.93%  '
$ws.Range("E33").Value = '  -0.89%  '
$ws.Range("E34").Value = '  -3.35%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("E36").Value = '  +2.26%  '
$ws.Range("E37").Value = '  -4.61%  '
$ws.Range("E38").Value = '  +0.24%  '
$ws.Range("E39").Value = '  -4.37%  '
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("E41").Value = '  -3.43%  '
$ws.Range("E42").Value = '  -3.31%  '
$ws.Range("E43").Value = '  -7.78%  '
$ws.Range("E44").Value = '  -13.69%  '
$ws.Range("E45").Value = '  -10.00%  '
$ws.Range("E46").Value = '  -2.62%  '
$ws.Range("E47").Value = '  -6.43%  '
$ws.Range("E48").Value = '  -3.68%  '
$ws.Range("E49").Value = '  -3.84%  '
$ws.Range("E50").Value = '  -5.45%  '
$ws.Range("E51").Value = '  -2.55%  '

# Restore original (default) cell style now that values are entered as text,
# so no stray style/number-format attribute is left on the cells.
$ws.Range("D2:E51").Style = "Normal"
